$d = $word.ActiveDocument

# Update the date heading at the top of the document
$d.Content.Find.Execute("2024-10-27 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-28 Monday", 2) | Out-Null

# Update each arithmetic expression cell in the table, in row/column order
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "0+99="
$t.Cell(1, 2).Range.Text = "73-67="
$t.Cell(1, 3).Range.Text = "11+13="
$t.Cell(1, 4).Range.Text = "78-40="
$t.Cell(1, 5).Range.Text = "19+14="

$t.Cell(2, 1).Range.Text = "64-53="
$t.Cell(2, 2).Range.Text = "57+12="
$t.Cell(2, 3).Range.Text = "26-26="
$t.Cell(2, 4).Range.Text = "36+59="
$t.Cell(2, 5).Range.Text = "6+59="

$t.Cell(3, 1).Range.Text = "28+22="
$t.Cell(3, 2).Range.Text = "21+55="
$t.Cell(3, 3).Range.Text = "79-38="
$t.Cell(3, 4).Range.Text = "23+9="
$t.Cell(3, 5).Range.Text = "28+3="

$t.Cell(4, 1).Range.Text = "56-17="
$t.Cell(4, 2).Range.Text = "17-14="
$t.Cell(4, 3).Range.Text = "71+15="
$t.Cell(4, 4).Range.Text = "56-47="
$t.Cell(4, 5).Range.Text = "14+16="

$t.Cell(5, 1).Range.Text = "15+32="
$t.Cell(5, 2).Range.Text = "55-34="
$t.Cell(5, 3).Range.Text = "29+21="
$t.Cell(5, 4).Range.Text = "41-40="
$t.Cell(5, 5).Range.Text = "10+88="

$t.Cell(6, 1).Range.Text = "92-89="
$t.Cell(6, 2).Range.Text = "25+62="
$t.Cell(6, 3).Range.Text = "23+46="
$t.Cell(6, 4).Range.Text = "76+3="
$t.Cell(6, 5).Range.Text = "96-82="

$t.Cell(7, 1).Range.Text = "95-74="
$t.Cell(7, 2).Range.Text = "58+4="
$t.Cell(7, 3).Range.Text = "1+93="
$t.Cell(7, 4).Range.Text = "95-86="
$t.Cell(7, 5).Range.Text = "12+61="

$t.Cell(8, 1).Range.Text = "72-4="
$t.Cell(8, 2).Range.Text = "8+46="
$t.Cell(8, 3).Range.Text = "0+53="
$t.Cell(8, 4).Range.Text = "60-55="
$t.Cell(8, 5).Range.Text = "43-39="

$t.Cell(9, 1).Range.Text = "74-52="
$t.Cell(9, 2).Range.Text = "68-32="
$t.Cell(9, 3).Range.Text = "82-70="
$t.Cell(9, 4).Range.Text = "85-43="
$t.Cell(9, 5).Range.Text = "58-51="

$t.Cell(10, 1).Range.Text = "5+9="
$t.Cell(10, 2).Range.Text = "23+13="
$t.Cell(10, 3).Range.Text = "32+13="
$t.Cell(10, 4).Range.Text = "45+43="
$t.Cell(10, 5).Range.Text = "30+33="

$t.Cell(11, 1).Range.Text = "34-33="
$t.Cell(11, 2).Range.Text = "11+21="
$t.Cell(11, 3).Range.Text = "51-34="
$t.Cell(11, 4).Range.Text = "77-21="
$t.Cell(11, 5).Range.Text = "16+70="

$t.Cell(12, 1).Range.Text = "76-73="
$t.Cell(12, 2).Range.Text = "43-38="
$t.Cell(12, 3).Range.Text = "9+45="
$t.Cell(12, 4).Range.Text = "17+46="
$t.Cell(12, 5).Range.Text = "25-23="

$t.Cell(13, 1).Range.Text = "19+61="
$t.Cell(13, 2).Range.Text = "9+54="
$t.Cell(13, 3).Range.Text = "64+22="
$t.Cell(13, 4).Range.Text = "87-48="
$t.Cell(13, 5).Range.Text = "33-19="

$t.Cell(14, 1).Range.Text = "70-34="
$t.Cell(14, 2).Range.Text = "3+30="
$t.Cell(14, 3).Range.Text = "63-31="
$t.Cell(14, 4).Range.Text = "40-29="
$t.Cell(14, 5).Range.Text = "0+11="

$t.Cell(15, 1).Range.Text = "24+25="
$t.Cell(15, 2).Range.Text = "24+27="
$t.Cell(15, 3).Range.Text = "2+12="
$t.Cell(15, 4).Range.Text = "96-11="
$t.Cell(15, 5).Range.Text = "20+31="

$t.Cell(16, 1).Range.Text = "63+14="
$t.Cell(16, 2).Range.Text = "44-15="
$t.Cell(16, 3).Range.Text = "95-59="
$t.Cell(16, 4).Range.Text = "36+28="
$t.Cell(16, 5).Range.Text = "3+18="

$t.Cell(17, 1).Range.Text = "93-75="
$t.Cell(17, 2).Range.Text = "20+2="
$t.Cell(17, 3).Range.Text = "86-62="
$t.Cell(17, 4).Range.Text = "47+45="
$t.Cell(17, 5).Range.Text = "26-8="

$t.Cell(18, 1).Range.Text = "88+2="
$t.Cell(18, 2).Range.Text = "85-50="
$t.Cell(18, 3).Range.Text = "81-65="
$t.Cell(18, 4).Range.Text = "78-60="
$t.Cell(18, 5).Range.Text = "52+45="

$t.Cell(19, 1).Range.Text = "55-3="
$t.Cell(19, 2).Range.Text = "15+61="
$t.Cell(19, 3).Range.Text = "57+10="
$t.Cell(19, 4).Range.Text = "21-19="
$t.Cell(19, 5).Range.Text = "58-29="

$t.Cell(20, 1).Range.Text = "56+0="
$t.Cell(20, 2).Range.Text = "47+9="
$t.Cell(20, 3).Range.Text = "7+56="
$t.Cell(20, 4).Range.Text = "76-43="
$t.Cell(20, 5).Range.Text = "67+14="
